$wb = $excel.ActiveWorkbook
try {
  $s = $wb.Styles.Add("TempStyle")
  Write-Host "Added temp style"
} catch {
  Write-Host "Error Add: $_"
}
